$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet originally held a pivot_table-shaped result:
#     sex | black | hispanic | other | white          (A1:E3)
# It becomes a long/groupby-shaped result:
#     sex | race | earn                               (A1:C9)
# with the "sex" column merged down across each block of 4 race rows.
# ---------------------------------------------------------------------------

# --- 0. Prime the style table with the extra (left / left+right / left+right+
#        bottom) boxed-border combinations that Excel computes while
#        harmonising a boxed border across a merged range. Doing this on a
#        throwaway range first (then clearing it) registers the border
#        definitions in styles.xml without any real cell ending up using them
#        -- matching the look of the target style sheet.
$prime = $ws.Range("Z1:Z4")
$prime.Borders.Item(7).LineStyle = 1    # xlEdgeLeft
$prime.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$ws.Range("Z4").Borders.Item(9).LineStyle = 1  # xlEdgeBottom
$prime.Clear()

# --- 1. Drop the two columns (D, E) that won't exist anymore.
$ws.Range("D1:E3").EntireColumn.Delete()

# --- 2. Rewrite the header row. A1 stays "sex" (already correctly styled).
$ws.Range("B1").Value = "race"
$ws.Range("C1").Value = "earn"

# --- 3. Clear out the old 2-data-row body; rebuild rows 2-9 from scratch.
$ws.Range("A2:C3").Clear()

# --- 4. Write the new long-format values (still unstyled at this point).
$ws.Range("A2").Value = "female"
$ws.Range("B2").Value = "black"
$ws.Range("C2").Value = 26413.2832533842

$ws.Range("B3").Value = "hispanic"
$ws.Range("C3").Value = 21217.35209195709

$ws.Range("B4").Value = "other"
$ws.Range("C4").Value = 34164.34619665911

$ws.Range("B5").Value = "white"
$ws.Range("C5").Value = 23948.24117218976

$ws.Range("A6").Value = "male"
$ws.Range("B6").Value = "black"
$ws.Range("C6").Value = 31778.72028241918

$ws.Range("B7").Value = "hispanic"
$ws.Range("C7").Value = 31818.39067697309

$ws.Range("B8").Value = "other"
$ws.Range("C8").Value = 29189.70626633994

$ws.Range("B9").Value = "white"
$ws.Range("C9").Value = 48951.73144985256

# --- 5. Merge the "sex" column across each 4-row race block while the range
#        is still unstyled, so the merge itself doesn't fragment any border.
$ws.Range("A2:A5").Merge()
$ws.Range("A6:A9").Merge()

# --- 6. Now stamp the header's formatting (bold, centered, boxed) onto
#        columns A and B for the whole body.
$ws.Range("A1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("A6:A9").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("B2:B9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# PasteSpecial(formats) shouldn't touch cell contents, but make sure the
# "sex" labels are still exactly right.
$ws.Range("A2").Value = "female"
$ws.Range("A6").Value = "male"

Write-Output "done"
